# Trade #104 closed at 2026-02-17 09:18:25 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) "Summary" sheet - refresh aggregate stats after new trade close
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.1    # Current Capital
$summary.Range("B4").Value = 0.11      # Total P&L $
$summary.Range("B5").Value = 0.02      # Total P&L %
$summary.Range("B6").Value = 104       # Total Trades
$summary.Range("B7").Value = 44        # Winning Trades
$summary.Range("B9").Value = 42.31     # Win Rate %

# ---------------------------------------------------------------
# 2) "Strategy Status" sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.1      # Capital
$status.Range("D4").Value = 104        # Trades
$status.Range("E4").Value = 0.11       # P&L $
$status.Range("F4").Value = 0.1        # P&L %
$status.Range("G4").Value = 42.31      # Win Rate %

# ---------------------------------------------------------------
# 3) "All Trades" sheet - append trade #104 as new row 105
# ---------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("A105").Value = 104
# Leading apostrophe forces text (prevents auto-conversion of the
# ISO-looking date string into a real date serial); resetting the
# style back to Normal afterwards drops the resulting quote-prefix
# formatting so the cell keeps the sheet's default (unstyled) look.
$allTrades.Range("B105").Value = "'2026-02-17"
$allTrades.Range("B105").Style = "Normal"
$allTrades.Range("C105").Value = "09:18:19"
$allTrades.Range("D105").Value = "MarketMaking"
$allTrades.Range("E105").Value = "DOWN"
$allTrades.Range("F105").Value = 0.72
$allTrades.Range("G105").Value = 0.77
$allTrades.Range("H105").Value = "CLOSED"
$allTrades.Range("I105").Value = 6.9444
$allTrades.Range("J105").Value = 0.05
$allTrades.Range("K105").Value = 100.1
$allTrades.Range("L105").Value = 0
$allTrades.Range("M105").Value = 0
$allTrades.Range("N105").Value = 0.6
$allTrades.Range("O105").Value = "Normal spread capture: 19600 bps"
$allTrades.Range("P105").Value = "early_exit"
$allTrades.Range("Q105").Value = 0.13

# ---------------------------------------------------------------
# 4) "MarketMaking" sheet - same new trade row, mirrored log
# ---------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Range("A105").Value = 104
$marketMaking.Range("B105").Value = "'2026-02-17"
$marketMaking.Range("B105").Style = "Normal"
$marketMaking.Range("C105").Value = "09:18:19"
$marketMaking.Range("D105").Value = "MarketMaking"
$marketMaking.Range("E105").Value = "DOWN"
$marketMaking.Range("F105").Value = 0.72
$marketMaking.Range("G105").Value = 0.77
$marketMaking.Range("H105").Value = "CLOSED"
$marketMaking.Range("I105").Value = 6.9444
$marketMaking.Range("J105").Value = 0.05
$marketMaking.Range("K105").Value = 100.1
$marketMaking.Range("L105").Value = 0
$marketMaking.Range("M105").Value = 0
$marketMaking.Range("N105").Value = 0.6
$marketMaking.Range("O105").Value = "Normal spread capture: 19600 bps"
$marketMaking.Range("P105").Value = "early_exit"
$marketMaking.Range("Q105").Value = 0.13
